$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Espinaca" (spinach)
# series for "Feria Lagunitas de Puerto Montt". It belongs right above the
# existing row 79, so push that row (and everything below it, through the
# old last row 102) down by one -- exactly like Excel's "Insert Sheet Rows"
# on the row-79 header.
$ws.Rows.Item(79).Insert()

# Populate the freshly inserted row 79 with the new observation.
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 45275
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112012
$ws.Range("G79").Value = "Espinaca"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 40
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = 15000
$ws.Range("N79").Value = "`$/cuna 10 kilos"
$ws.Range("O79").Value = "Región Metropolitana"
$ws.Range("P79").Value = 1500
$ws.Range("Q79").Value = 10
$ws.Range("R79").Value = "Hortaliza"
